$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): updates to column F ("想去人数")
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 158
$ws1.Range("F5").Value = 4706
$ws1.Range("F8").Value = 520
$ws1.Range("F13").Value = 2836
$ws1.Range("F15").Value = 94
$ws1.Range("F18").Value = 2373
$ws1.Range("F19").Value = 105
$ws1.Range("F25").Value = 234
$ws1.Range("F26").Value = 43

# Sheet "全部类型" (sheet4): same underlying event rows, offset by one row
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 158
$ws4.Range("F6").Value = 4706
$ws4.Range("F9").Value = 520
$ws4.Range("F14").Value = 2836
$ws4.Range("F16").Value = 94
$ws4.Range("F19").Value = 2373
$ws4.Range("F20").Value = 105
$ws4.Range("F26").Value = 234
$ws4.Range("F27").Value = 43
